$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the explicit per-cell / per-column style (font) that used to be
# applied across the whole sheet; headers go back to the default style.
$ws.Cells.ClearFormats()

# Insert a new column before the last one; this shifts the existing
# "Validation CFVU" header from column I to column J.
$ws.Columns.Item(9).Insert()

# Set the new header value in the newly inserted column I.
$ws.Cells.Item(1, 9).Value = "Site formation"

# Column E grew wider (content used to autofit against longer values).
$ws.Columns.Item(5).ColumnWidth = 18.166666666666664

# New columns F:J get their own explicit widths.
$ws.Columns.Item(6).ColumnWidth = 18.333333333333336
$ws.Columns.Item(7).ColumnWidth = 9.333333333333332
$ws.Columns.Item(8).ColumnWidth = 15.666666666666668
$ws.Columns.Item(9).ColumnWidth = 12.166666666666668
$ws.Columns.Item(10).ColumnWidth = 13.666666666666668

# Update the active selection to match the new authored state.
$ws.Range("K9").Select()

$wb.Save()
